$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.552.92"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "2.660.77"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.77"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.01"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("E8").Value = "  +5.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.124"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.85"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -1.77%  "

$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.31"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -2.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000196"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -3.87%  "

$ws.Range("D15").Value = "3.134.35"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "65.381.44"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "2.646.74"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.64"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("E19").Value = "  -1.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +2.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.88"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -2.00%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("E23").Value = "  -1.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000110"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +4.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.71"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.62"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -4.64%  "

$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "533.50"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("E32").Value = "  -3.94%  "

$ws.Range("E33").Value = "  -1.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +2.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.36"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.08%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "159.11"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -2.39%  "

$ws.Range("E40").Value = "  -2.96%  "

$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "164.30"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.18%  "

$ws.Range("E44").Value = "  -1.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0607"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("E46").Value = "  -1.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.92"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0259"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.86%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.640"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -1.88%  "

$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("E51").Value = "  +2.15%  "
